# Update the cryptos list on the worksheet with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to be stored as text, even when
# the string content looks like a number (e.g. "112.61"), matching the
# original sheet where every Price/Volume cell is a plain text value.
function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $c = $ws.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.116.43"
Set-TextValue "E2" "  -1.60%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.259.36"
Set-TextValue "E3" "  -1.80%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.02%  "

# Row 5 - Solana
Set-TextValue "D5" "112.61"
Set-TextValue "E5" "  +3.20%  "

# Row 6 - BNB
Set-TextValue "D6" "263.88"
Set-TextValue "E6" "  -2.92%  "

# Row 7 - XRP
Set-TextValue "D7" "0.616"
Set-TextValue "E7" "  -1.65%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.17%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.598"
Set-TextValue "E9" "  -3.00%  "

# Row 10 - Avalanche
Set-TextValue "D10" "47.63"
Set-TextValue "E10" "  +0.91%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0921"
Set-TextValue "E11" "  -1.84%  "

# Row 12 - Polkadot
Set-TextValue "D12" "8.68"
Set-TextValue "E12" "  +3.04%  "

# Row 13 - TRON
Set-TextValue "E13" "  -0.77%  "

# Row 14 - Chainlink
Set-TextValue "D14" "15.37"
Set-TextValue "E14" "  -2.37%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.596.01"
Set-TextValue "E15" "  -1.82%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.853"
Set-TextValue "E16" "  -0.77%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.251.47"
Set-TextValue "E17" "  -1.95%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.991.83"
Set-TextValue "E18" "  -1.82%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0000107"
Set-TextValue "E19" "  -3.73%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.90"
Set-TextValue "E20" "  +9.48%  "

# Row 21 - Litecoin
Set-TextValue "D21" "70.83"
Set-TextValue "E21" "  -1.99%  "

# Row 22 - ImmutableX
Set-TextValue "D22" "2.40"
Set-TextValue "E22" "  -3.50%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "229.56"
Set-TextValue "E23" "  -2.01%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "9.64"
Set-TextValue "E24" "  +4.53%  "

# Row 25
Set-TextValue "D25" "2.82"
Set-TextValue "E25" "  -5.04%  "

# Row 26
Set-TextValue "E26" "  -0.03%  "

# Row 27
Set-TextValue "D27" "11.24"
Set-TextValue "E27" "  -1.24%  "

# Row 28 - LEO
Set-TextValue "D28" "3.86"
Set-TextValue "E28" "  -1.81%  "

# Row 29 - InjectiveProtocol
Set-TextValue "D29" "41.10"
Set-TextValue "E29" "  -0.21%  "

# Row 30 - WEMIXToken
Set-TextValue "E30" "  -2.29%  "

# Row 31 - Toncoin
Set-TextValue "E31" "  -1.45%  "

# Row 32 - Monero
Set-TextValue "D32" "171.30"
Set-TextValue "E32" "  -3.75%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "21.23"
Set-TextValue "E33" "  -3.24%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0895"
Set-TextValue "E34" "  -2.19%  "

# Row 35 - Filecoin
Set-TextValue "D35" "5.57"
Set-TextValue "E35" "  -0.80%  "

# Row 36 - Stellar
Set-TextValue "D36" "0.126"
Set-TextValue "E36" "  -0.65%  "

# Row 37 - RenderToken
Set-TextValue "D37" "4.61"
Set-TextValue "E37" "  -4.97%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.0349"
Set-TextValue "E38" "  -2.76%  "

# Row 39 - was Kaspa, now NEARProtocol (rows 39/40 content swapped)
Set-TextValue "B39" "NEARProtocol"
Set-TextValue "C39" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "3.74"
Set-TextValue "E39" "  +0.26%  "

# Row 40 - was NEARProtocol, now Kaspa
Set-TextValue "B40" "Kaspa"
Set-TextValue "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D40" "0.103"
Set-TextValue "E40" "  -9.83%  "

# Row 41 - Celestia
Set-TextValue "D41" "14.01"
Set-TextValue "E41" "  +15.13%  "

# Row 42 - MultiversX
Set-TextValue "D42" "74.20"
Set-TextValue "E42" "  +10.14%  "

# Row 43 - LidoDAOToken
Set-TextValue "D43" "2.39"
Set-TextValue "E43" "  +2.19%  "

# Row 44 - Algorand
Set-TextValue "D44" "0.233"
Set-TextValue "E44" "  -1.63%  "

# Row 45 - THORChain
Set-TextValue "D45" "6.09"
Set-TextValue "E45" "  +10.91%  "

# Row 46 - FirstDigitalUSD
Set-TextValue "E46" "  -0.20%  "

# Row 47 - ARBITRUM
Set-TextValue "D47" "1.36"
Set-TextValue "E47" "  -2.47%  "

# Row 48 - FraxShare
Set-TextValue "D48" "8.53"
Set-TextValue "E48" "  -3.05%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0986"
Set-TextValue "E49" "  -3.12%  "

# Row 50 - was TrustWalletToken, now Aave (rows 50/51 content swapped)
Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "99.79"
Set-TextValue "E50" "  +0.29%  "

# Row 51 - was Aave, now TrustWalletToken
Set-TextValue "B51" "TrustWalletToken"
Set-TextValue "C51" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D51" "1.23"
Set-TextValue "E51" "  -0.48%  "
